$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stimulus")

# Fill in row 5 with new test plan data
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = "Basic transfers"
$ws.Range("D5").Value = "Multiple Write + Multiple Read txns with hsize, haddr same for each set of txns, to test on hsize handling."

# Update selection to C6
$ws.Range("C6").Select()
